$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.197866439819336
$ws.Range("B1").Value = 2.447990417480469
$ws.Range("C1").Value = 1.811211585998535
$ws.Range("D1").Value = 1.688495278358459
$ws.Range("E1").Value = 1.61900782585144
